# Generate Report for Handback
#
# This script updates the localization-status workbook to reflect that the
# handback (target-language files coming back from translation) has
# completed and is back in sync with en-US.
#
#  - "Ready for handoff" status becomes "Handed back: in sync with en-US"
#    on every sheet (Overview + each language tab).
#  - The per-language detail tabs (zh-cn, de-de) get their "Latest Target
#    File" / "Latest Handback File" / "Latest Handback DateTime" columns
#    populated for both rows, including a hyperlink for the target file.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

$mdUrlBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6a16ac6059858f59bdfe9dd99051b8a13a1ea242/e2e/"

# ---------------------------------------------------------------------
# Overview sheet: replace the status text in columns E (zh-cn) and F (de-de)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------
# Per-language detail sheets: zh-cn finished handback at 16:31:42,
# de-de finished handback at 16:31:50.
# ---------------------------------------------------------------------
$languages = @(
    @{ Sheet = "zh-cn"; HandbackTime = "2016-08-23 16:31:42" },
    @{ Sheet = "de-de"; HandbackTime = "2016-08-23 16:31:50" }
)

$fileRows = @(
    @{ Row = 2; Slug = "0cef30ca-af02-42a2-ae69-598dcc329269"; Hash = "b892b91837adc3b147428ad887877575e06bf0d0" },
    @{ Row = 3; Slug = "9ade4afb-96ea-4e71-921a-8bccd71c9f25"; Hash = "fd38582eafb986306c4d859b1e11fc3c7d82154d" }
)

foreach ($lang in $languages) {
    $ws = $wb.Worksheets.Item($lang.Sheet)

    foreach ($f in $fileRows) {
        $row = $f.Row
        $mdName = "$($f.Slug).md"

        # Status column (C)
        $ws.Range("C$row").Value = $newStatus

        # Latest Target File (I) / Latest Handback File (J) / Latest Handback DateTime (K)
        $xlfName = "$($f.Slug).$($f.Hash).$($lang.Sheet).xlf"
        $ws.Range("I$row").Value = $mdName
        $ws.Range("J$row").Value = $xlfName
        $ws.Range("K$row").Value = $lang.HandbackTime
    }

    $ws.Columns.Item(3).ColumnWidth = 29.166666666666668
    $ws.Columns.Item(9).ColumnWidth = 39.166666666666664
    $ws.Columns.Item(10).ColumnWidth = 39.166666666666664

    # Re-create hyperlinks in final cell order (A2, I2, A3, I3) so the
    # relationship ids line up the same way Excel would renumber them.
    $ws.Hyperlinks.Delete()
    foreach ($f in $fileRows) {
        $row = $f.Row
        $mdName = "$($f.Slug).md"
        $url = "$mdUrlBase$mdName"
        $ws.Hyperlinks.Add($ws.Range("A$row"), $url, $null, $null, $mdName) | Out-Null
        $ws.Hyperlinks.Add($ws.Range("I$row"), $url, $null, $null, $mdName) | Out-Null
    }

    # Match the hyperlink text style ("HyperLink") on the newly-linked cells.
    $ws.Range("I2").Style = "HyperLink"
    $ws.Range("I3").Style = "HyperLink"
}
